$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "registro_docente" table but keep its data/contents ---
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Unlist()
}

# --- Remove the two bogus trailing rows (old "ss" test row + old Sebastian row) ---
$ws.Rows.Item(3).Resize(2).Delete()

# --- Trim row 2 down to columns A:G (drop Oficina/Departamento/Ciudad=H,I,J) ---
$ws.Range("H2:J2").ClearContents()

# --- Replace the data row (was Maryem Ruiz) with Sebastian Palacio's record ---
$ws.Range("A2").Value = "Sebastián"
$ws.Range("B2").Value = "Palacio"
$ws.Range("C2").Value = 1000762620
$ws.Range("D2").Value = "sebasx200"

# E2 must be stored as TEXT "1234" (not a number). Build it via a formula and
# paste-special as values so it lands as a shared string with no extra style.
$ws.Range("E2").Formula = "=""1234"""
$ws.Range("E2").Copy()
$ws.Range("E2").PasteSpecial(-4163)

$ws.Range("G2").Value = "No tiene"

# --- Correo (F2): fix the text + repoint the hyperlink ---
$ws.Range("F2").Hyperlinks.Delete()
$ws.Range("F2").Value = "sebastian_palacio23231@elpoli.edu.co"
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:sebastian_palacio23231@elpoli.edu.co")
$ws.Range("F2").Style = "Hipervínculo"

# --- Header row: drop the centered style so headers use the default style ---
$ws.Range("A1:J1").ClearFormats()

$excel.CutCopyMode = 0
